# Re-run of the last logged test-case block for each sheet (log was
# re-appended for each sheet — "path changed" re-run of the same suite).
#
# For every sheet we duplicate the trailing block of log rows (the same
# rows that were written for the previous run of that sheet's test) onto
# the end of the sheet, reusing the existing shared-string values so no
# new unique strings are introduced.

$wb = $excel.ActiveWorkbook

function Append-Block {
    param(
        [string]$SheetName,
        [int]$BlockSize
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Find the last used row in column A.
    $lastRow = $ws.Cells.Item(1048576, 1).End(-4162).Row

    $startOfBlock = $lastRow - $BlockSize + 1

    for ($i = 0; $i -lt $BlockSize; $i++) {
        $srcRow = $startOfBlock + $i
        $destRow = $lastRow + 1 + $i
        $val = $ws.Cells.Item($srcRow, 1).Value2
        $ws.Cells.Item($destRow, 1).Value = $val
    }
}

Append-Block "TestCase1_HospitalFilter" 17
Append-Block "TestCase0_MaxRatingFinder" 5
Append-Block "TestCase2_TopCities" 9
Append-Block "TestCase3_InvalidForm" 7
